# Reorder the tied_teams list entries (column O) for rows where three (or four)
# teams tied in European/FIFA World Cup qualification groups, per the commit:
#   "evaluate suspense for three teams tied in European Championships"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40-52: ['Ireland', 'Costa Rica'] -> ['Costa Rica', 'Ireland']
foreach ($r in 40..52) {
    $ws.Range("O$r").Value = "['Costa Rica', 'Ireland']"
}

# Rows 53-59: ['Colombia', 'Ireland', 'Argentina', 'Costa Rica'] -> ['Costa Rica', 'Colombia', 'Argentina', 'Ireland']
foreach ($r in 53..59) {
    $ws.Range("O$r").Value = "['Costa Rica', 'Colombia', 'Argentina', 'Ireland']"
}

# Rows 63-73: ['Scotland', 'Colombia', 'Argentina', 'Austria'] -> ['Colombia', 'Scotland', 'Austria', 'Argentina']
foreach ($r in 63..73) {
    $ws.Range("O$r").Value = "['Colombia', 'Scotland', 'Austria', 'Argentina']"
}

# Row 78: ['South Korea', 'Netherlands'] -> ['Netherlands', 'South Korea']
$ws.Range("O78").Value = "['Netherlands', 'South Korea']"
